$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the new font/style (Calibri 10, centered) to the data blocks
# that previously used style index 2 (rows 7-15 and 17-25, columns C:M).
$rng1 = $ws.Range("C7:M15")
$rng1.Font.Name = "Calibri"
$rng1.Font.Size = 10
$rng1.HorizontalAlignment = -4108
$rng1.VerticalAlignment = -4108

$rng2 = $ws.Range("C17:M25")
$rng2.Font.Name = "Calibri"
$rng2.Font.Size = 10
$rng2.HorizontalAlignment = -4108
$rng2.VerticalAlignment = -4108

# Correct values for the TSL-3000-60-1-230-IP67 row in the height table (row 13)
$ws.Range("L13").Value = 515
$ws.Range("M13").Value = 515

# Correct values for the TSL-3000-60-1-230-IP67 row in the weight table (row 23)
$ws.Range("L23").Value = 55
$ws.Range("M23").Value = 80

# Update the saved selection/active cell
$ws.Range("T28").Select()
